# ProjectObjects.xlsx -- "Class" sheet restructure.
#
# The "Class" sheet's lower block (originally rows 10-13, columns D..I) is
# renumbered to rows 12-15 (one column to the left, C..H) and a brand-new
# row 11 is inserted above it holding a single new label, "Simple Screen",
# in column C.
#
#   old D10:I13                       new C11:H15
#   ----------------------------      -----------------------------------
#   D10 Gold                     ->   C12 Gold
#   D11 Hero[]  E11 Name          ->   C13 Hero[]  D13 Name
#   F11 DamagePerClick            ->   E13 DamagePerClick
#   G11 DamagePerSecond           ->   F13 DamagePerSecond
#   H11 Level   I11 UpgradeGold   ->   G13 Level   H13 UpgradeGold
#   D12 Monster E12 Name F12 HP   ->   C14 Monster D14 Name  E14 HP
#   D13 Level E13 Number          ->   C15 Level   D15 Number
#   F13 Name  G13 Monster Killed  ->   E15 Name    F15 Monster Killed
#                                      C11 Simple Screen   (new)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Class")

# 1) Insert two fresh blank rows above the old row 10. This pushes the
#    existing rows 10-13 down to 12-15 (row 10 stays empty/unused and row
#    11 becomes the new row that will carry the "Simple Screen" label).
$ws.Range("A10:A11").EntireRow.Insert()

# Give the newly-inserted row 11 the same height/formatting the other
# rows already use.
$ws.Rows.Item(11).RowHeight = 29.25

# 2) The data that used to live in columns D..I now needs to shift one
#    column to the left (into C..H) on each of the relocated rows.
foreach ($r in 12, 13, 14, 15) {
    for ($col = 4; $col -le 9; $col++) {
        $srcCell = $ws.Cells.Item($r, $col)
        $dstCell = $ws.Cells.Item($r, $col - 1)
        $dstCell.Value = $srcCell.Value()
        $srcCell.Clear()
    }
}

# 3) Populate the new row with its label.
$ws.Range("C11").Value = "Simple Screen"

# 4) Match the author's resulting selection/view.
$ws.Range("C12").Select()
